$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and E (Volume(1h)) cell text updates per the commit diff.
# Cells whose new text parses as a plain number are pinned to Text format (NumberFormat "@")
# first, matching the original inlineStr/text storage, so Excel does not silently coerce
# strings like "1.00" or "144.70" into numeric 1 / 144.7 and drop the literal formatting.

$ws.Range("D2").Value = "26.280.29"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "1.594.26"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.68"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("E6").Value = "  +1.08%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.38"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "1.817.84"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "1.587.76"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.42"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "26.276.06"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.44"
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "213.76"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.02"
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("E24").Value = "  -2.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.70"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.20"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D34").Value = "1.337.94"
$ws.Range("E34").Value = "  +4.85%  "
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.592"
$ws.Range("E37").Value = "  -3.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0167"
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("E40").Value = "  -11.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.71"
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("D46").Value = "1.729.58"
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.82"
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("E48").Value = "  -3.98%  "
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0977"
$ws.Range("E50").Value = "  -2.79%  "

Write-Host "Applied cryptos list update."
